# Add a new "2022-Q4" detail sheet (between "总计" and "2021-Q1") and record
# its summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item(1)
$ws2021 = $wb.Worksheets.Item(2)

function Set-TextValue($range, [string]$text) {
    # Force a numeric-looking string (fund codes with leading zeros, "6.19",
    # "94.10", ...) to be stored as text instead of being coerced to a
    # number, without leaving a stray NumberFormat/quotePrefix style behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for 2022-Q4, pushing 2021-Q1 to row 3
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

# Carry the "index" column's number style (s=2) down onto the newly
# inserted row, then clear the stray formatting Insert() left on B:D.
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$zongji.Range("B2:D2").ClearFormats()

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 6
$zongji.Range("D2").Value = 0.53

$zongji.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q4" worksheet right before the existing "2021-Q1"
#    sheet (so the tab order becomes 总计, 2022-Q4, 2021-Q1).
# ---------------------------------------------------------------------
$newWs = $wb.Worksheets.Add($ws2021)
$newWs.Name = "2022-Q4"

# Borrow the "总计" sheet's header style (s=2) for the header row and the
# index column, matching how this workbook styles its detail sheets.
$zongji.Range("B1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$newWs.Range("A2:A7").PasteSpecial(-4122)

$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

$data = @(
    @(0, "011184", "东方阿尔法招阳混合A",               "6.19", "94.10", "6.86", "0.4246", 6),
    @(1, "011185", "东方阿尔法招阳混合C",               "0.88", "94.10", "6.86", "0.0604", 6),
    @(2, "350005", "天治中国制造2025灵活配置混合",       "0.24", "94.00", "6.16", "0.0148", 3),
    @(3, "350009", "天治研究驱动混合A",                 "0.29", "92.96", "3.90", "0.0113", 1),
    @(4, "350002", "天治低碳经济灵活配置混合",           "0.65", "85.45", "1.43", "0.0093", 6),
    @(5, "002043", "天治研究驱动混合C",                 "0.17", "92.96", "3.90", "0.0066", 1)
)

$row = 2
foreach ($rec in $data) {
    $newWs.Range("A$row").Value = $rec[0]
    Set-TextValue $newWs.Range("B$row") $rec[1]
    Set-TextValue $newWs.Range("C$row") $rec[2]
    Set-TextValue $newWs.Range("D$row") $rec[3]
    Set-TextValue $newWs.Range("E$row") $rec[4]
    Set-TextValue $newWs.Range("F$row") $rec[5]
    Set-TextValue $newWs.Range("G$row") $rec[6]
    $newWs.Range("H$row").Value = $rec[7]
    $row = $row + 1
}

# "2021-Q1" was the active/selected tab before this edit; keep it that way
# (it now lives in a relocated sheet part, but should stay the active tab).
$wb.Worksheets.Item("2021-Q1").Activate()
